$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.978.02'
$ws.Range("E2").Value = '  -0.46%  '

$ws.Range("D3").Value = '1.826.45'
$ws.Range("E3").Value = '  +0.19%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.007'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.26%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.97'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.09%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4615'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.36%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3704'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.78%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07335'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.52%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8744'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.47%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07965'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.03%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '19.79'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.66%  '

$ws.Range("D13").Value = '1.792.23'
$ws.Range("E13").Value = '  -2.19%  '

$ws.Range("E14").Value = '  -0.15%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.555'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.23%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.29'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.42%  '

$ws.Range("E17").Value = '  -0.12%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008865'
$ws.Range("D18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.007'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.14%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.81'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.35%  '

$ws.Range("D21").Value = '27.246.95'
$ws.Range("E21").Value = '  -0.58%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.107'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.77%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.54'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.12%  '

$ws.Range("D24").Value = '2.115.67'
$ws.Range("E24").Value = '  +1.05%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.16'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.95%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.847'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.51%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.43'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.93%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.040'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.82%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.137'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.53%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '115.39'
$ws.Range("D30").Style = "Normal"

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08904'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.12%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.967'
$ws.Range("D32").Style = "Normal"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7287'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.80%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.432'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.50%  '

$ws.Range("E35").Value = '  -0.93%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.478'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.61%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01953'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.63%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.070'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.30%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05228'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.57%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.946'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.34%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.085'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.03%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5156'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.06%  '

$ws.Range("E43").Value = '  -0.43%  '

$ws.Range("E44").Value = '  -1.09%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4843'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.92%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.25'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.75%  '

$ws.Range("E47").Value = '  -0.34%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '102.65'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.20%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.631'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.17%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06196'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.96%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '64.66'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.09%  '

